$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.011636972427368
$ws.Range("B1").Value = 1.645598530769348
$ws.Range("C1").Value = 6.794440269470215
$ws.Range("D1").Value = 2.687289476394653
$ws.Range("E1").Value = 1.483729839324951
